$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("R1")
$dst = $ws.Range("S1")

$dst.Value = "NOMBRE(s) ARCHIVO EVIDENCIA"

$src.Copy()
$dst.PasteSpecial(-4122)

$dst.Borders.Item(8).LineStyle = -4142
$dst.Borders.Item(9).LineStyle = -4142
